# "add staff id and phone number to line manager"
#
# The sheet currently has 7 data columns (A:G). Column G holds each
# employee's second phone number. We insert a new column H for that
# phone number (so the "line manager" columns read: staff id, first
# name, last name, middle name, email, phone 1, phone 2) and leave G
# blank (still formatted, just empty) in the process, matching the
# target layout where the old G content slides over into the new H
# column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh column at H so it inherits the neighbouring column's
# formatting (style/width) instead of picking up the sheet default —
# mirrors how the existing G column already looks. Then drop the extra
# column that the insert pushed off the tracked-columns range so the
# sheet keeps its original 256-column width definition.
$ws.Columns.Item(8).Insert()
$ws.Columns.Item(257).Delete()

# Move each row's phone-number-2 value (and its style) from G to the
# newly created H, leaving G blank but still styled.
for ($r = 1; $r -le 10; $r++) {
    $src = $ws.Cells.Item($r, 7)
    $dst = $ws.Cells.Item($r, 8)
    $src.Copy($dst)
    $src.ClearContents()
}
